$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G10").Value = 2.5
$ws.Range("H10").Value = 2.7
$ws.Range("J10").Value = 3.4
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("AG10").Value = 7.5
$ws.Range("AJ10").Value = 34
$ws.Range("AL10").Value = 41
$ws.Range("AX10").Value = 19
$ws.Range("BA10").Value = 101
